$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-16 from 45207 to 45208
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 45208
}

# Update hyperlink formulas in row 2 (S2, T2, V2, W2, X2, Y2):
# replace "Logging_LINDESBERG" with "Logging_1885" in the URL path
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/artfynd/A 34293-2023.xlsx", "A 34293-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/kartor/A 34293-2023.png", "A 34293-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/klagomål/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/klagomålsmail/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/tillsyn/A 34293-2023.docx", "A 34293-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/tillsynsmail/A 34293-2023.docx", "A 34293-2023")'
